{"js": "const replacements = [\n  [\"94\u00d755=\", \"46\u00d796=\"],\n  [\"56\u00d724=\", \"88\u00d727=\"],\n  [\"65\u00d761=\", \"84\u00d774=\"],\n  [\"13\u00d770=\", \"79\u00d712=\"],\n  [\"81\u00d799=\", \"84\u00d762=\"],\n  [\"87\u00d748=\", \"19\u00d712=\"],\n  [\"31\u00d773=\", \"53\u00d796=\"],\n  [\"15\u00d798=\", \"17\u00d738=\"],\n  [\"48\u00d726=\", \"18\u00d777=\"],\n  [\"29\u00d741=\", \"98\u00d727=\"],\n  [\"95\u00d790=\", \"17\u00d757=\"],\n  [\"99\u00d787=\", \"79\u00d795=\"],\n  [\"18\u00d797=\", \"71\u00d720=\"],\n  [\"54\u00d715=\", \"96\u00d747=\"],\n  [\"95\u00d765=\", \"50\u00d769=\"],\n  [\"84\u00d775=\", \"77\u00d763=\"],\n  [\"14\u00d720=\", \"62\u00d725=\"],\n  [\"81\u00d797=\", \"88\u00d748=\"],\n  [\"27\u00d778=\", \"84\u00d799=\"],\n  [\"71\u00d780=\", \"22\u00d713=\"],\n  [\"30\u00d731=\", \"94\u00d763=\"],\n  [\"40\u00d788=\", \"33\u00d718=\"],\n  [\"49\u00d758=\", \"12\u00d770=\"],\n  [\"15\u00d799=\", \"22\u00d751=\"],\n  [\"39\u00d746=\", \"86\u00d784=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"94\u00d755=\", \"46\u00d796=\"),\n    @(\"56\u00d724=\", \"88\u00d727=\"),\n    @(\"65\u00d761=\", \"84\u00d774=\"),\n    @(\"13\u00d770=\", \"79\u00d712=\"),\n    @(\"81\u00d799=\", \"84\u00d762=\"),\n    @(\"87\u00d748=\", \"19\u00d712=\"),\n    @(\"31\u00d773=\", \"53\u00d796=\"),\n    @(\"15\u00d798=\", \"17\u00d738=\"),\n    @(\"48\u00d726=\", \"18\u00d777=\"),\n    @(\"29\u00d741=\", \"98\u00d727=\"),\n    @(\"95\u00d790=\", \"17\u00d757=\"),\n    @(\"99\u00d787=\", \"79\u00d795=\"),\n    @(\"18\u00d797=\", \"71\u00d720=\"),\n    @(\"54\u00d715=\", \"96\u00d747=\"),\n    @(\"95\u00d765=\", \"50\u00d769=\"),\n    @(\"84\u00d775=\", \"77\u00d763=\"),\n    @(\"14\u00d720=\", \"62\u00d725=\"),\n    @(\"81\u00d797=\", \"88\u00d748=\"),\n    @(\"27\u00d778=\", \"84\u00d799=\"),\n    @(\"71\u00d780=\", \"22\u00d713=\"),\n    @(\"30\u00d731=\", \"94\u00d763=\"),\n    @(\"40\u00d788=\", \"33\u00d718=\"),\n    @(\"49\u00d758=\", \"12\u00d770=\"),\n    @(\"15\u00d799=\", \"22\u00d751=\"),\n    @(\"39\u00d746=\", \"86\u00d784=\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Execute([ref]$oldText, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$newText, 2) | Out-Null\n}\n"}
